$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Update-Query {
    param(
        [string]$CellAddr
    )
    $cell = $ws.Range($CellAddr)
    $text = $cell.Value2
    $text = $text.Replace('df_participant prt ON std.id = prt."study.id"', 'df_participant prt ON std.study_id = prt."study.study_id"')
    $text = $text.Replace('df_diagnoses dgn ON prt.id = dgn."participant.id"', 'df_diagnoses dgn ON prt.participant_id = dgn."participant.participant_id"')
    $text = $text.Replace('df_treatments trt ON prt.id = trt."participant.id"', 'df_treatments trt ON prt.participant_id = trt."participant.participant_id"')
    $text = $text.Replace('df_treatment_resp trr ON prt.id = trr."participant.id"', 'df_treatment_resp trr ON prt.participant_id = trr."participant.participant_id"')
    $text = $text.Replace('df_survival srv ON prt.id = srv."participant.id"', 'df_survival srv ON prt.participant_id = srv."participant.participant_id"')
    $text = $text.Replace('df_reference_files rfs ON std.id = rfs."study.id"', 'df_reference_files rfs ON std.study_id = rfs."study.study_id"')
    $cell.Value2 = $text
}

# Update the 7 SQL query cells: the count query in C2, and the tab queries in B2:B7
Update-Query "C2"
Update-Query "B2"
Update-Query "B3"
Update-Query "B4"
Update-Query "B5"
Update-Query "B6"
Update-Query "B7"

# Update column C width (no longer best-fit, explicit width ~67.33)
$ws.Columns.Item(3).ColumnWidth = 66.5

# Update the active selection on the sheet
$ws.Range("B2").Select()
